$wb = $excel.ActiveWorkbook

# --- Update shared text on the "Score" sheet (also reflected on "Rank") ---
$wsScore = $wb.Worksheets.Item("Score")
$wsRank  = $wb.Worksheets.Item("Rank")

$wsScore.Range("B2").Value = "0.25:0.32:0.11:0.19:0.13"
$wsScore.Range("C2").Value = "0.25T, 0.32RR, 0.11Env, 0.19Econ, 0.13S"
$wsRank.Range("B2").Value = "0.25:0.32:0.11:0.19:0.13"
$wsRank.Range("C2").Value = "0.25T, 0.32RR, 0.11Env, 0.19Econ, 0.13S"

# --- Update the TOPSIS scores on the "Score" sheet ---
$wsScore.Range("D2").Value = 0.1682817982949932
$wsScore.Range("E2").Value = 0.8156312408875149
$wsScore.Range("F2").Value = 0.176540863984829

# --- Update the ranking on the "Rank" sheet ---
$wsRank.Range("D2").Value = 3
$wsRank.Range("E2").Value = 1
$wsRank.Range("F2").Value = 2
